$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target grid: header row (A1:K1) plus 12 FDA recall records (A2:K13).
# Headers move to UPPER_SNAKE_CASE; every free-text field is normalized to
# plain Title Case with punctuation collapsed to single spaces.
#
# $null marks a cell the source diff leaves untouched -- most notably
# A2:A13 ("1995-01-24"), which must stay skipped or this host's COM layer
# "helpfully" reinterprets the re-typed string as a date serial, same as
# real Excel would.
$reportGrid = @(
    @("RECALL_CLASSIFICATION_DATE", "PRODUCT_TYPE", "CLASSIFICATION", "RECALL_NUMBER", "PRODUCT", "RECALLING_FIRM", "MANUFACTURER", "RECALL_INITIALIZATION_DATE", "REASON", "VOLUME", "DISTRIBUTION"),
    @($null, $null, "II", $null, "Myochrysine Gold Sodium Thiomalate Injection 50 Mg Ml In 10 Ml Multi Dose Vials Used For Rheumatoid Arthritis Under The Merck And King Pharmaceuticals Label", $null, "Merck And Company Inc", "December 18 1995", "Product Does Not Meet Antimicrobial Effectiveness Test", "Approximately 90 000 Unit Packages Remained On Market At Time Of Recall Initiation", "Nationwide Spain Peru"),
    @($null, $null, "III", $null, "Sodium Iodide I 131 Solution For Oral Therapeutic Use In 100 Mci And 50 Mci Vials Used For The Treatment Of Hyperthyroidism And Selected Cases Of Carcinoma Of The Thyroid", "Cis Us", "Cis Bio International Subsidiary Of Compagnie Oris Industrie Sa Cedex France", "December 15 1995", "Some Packages Contained A Package Insert For Sodium Iodide I 131 Capsules", "71 100 Mci Vials And 55 50 Mci Vials Were Distributed", $null),
    @($null, $null, "II", $null, "I V Extension Sets With T Connectors A Abbott Extension Set With T Connector", $null, "Abbott Laboratories Laurinburg North Carolina Hg Lot Suffix Abbott Health Products Barceloneta Puerto Rico H 1 Lot Suffix", "December 8 1995", "Some Of The Extension Sets Leak And Disconnect Between The Male Fitment Of The Device And The Luer Fitment Of The I", "2 120 180 Units Were Distributed Firm Estimated That 135 000 150 000 Units Remained On Market At Time Of Recall Initiation", $null),
    @($null, $null, "II", $null, "System 97 A Helium Charged Mobile Intra Aortic Balloon Pump Iabp Catalog 0998 00 0104 X For Use As A Patient Aid During Pre Intra Or Post Operative Open Heart Surgery Use Within Patients Demonstrating Unstable Angina Use Within Patients With Left Main Artery Occlusion Or Poor Left Ventricle Function A Power Supply B Front End Board C Solenoid Driver Board", $null, "Datascope Corporation Paramus New Jersey", "August 31 1995", "The System May Fail Because Of The Following Component Failures Power Supply Failure Of A Zener Diode Designated At Vr 9 In The System Power Supply Which Can Result In A Power On Condition Which Cannot Be Switched Off By Normal Means", "A 259 Units B 554 Units C 153 Units Were Distributed", $null),
    @($null, $null, "II", $null, "3 M Sarns Brand Perfusion System 9000 Gas Flow System Used To Provide And Monitor Co 2 Flow To The Patient During Cardiopulmonary Bypass", $null, "Sarns 3 M Health Care Ann Arbor Michigan", "July 24 1995", "The Relay Controlling The Co 2 Flow Can Experience Random Failure And Latch In The On Position Causing Co 2 Gas To Continue To Flow After Completion Of The Pre Bypass Co 2 Flush Of The Perfusion Circuit", "Approximately 500 Units", $null),
    @($null, $null, "II", $null, "Duopulse And Unipulse Dental Laser System Used In Dental Surgery", "Manufacturer Fda Approved The Firm's Corrective Action Plan December 28", "Excel Quantronix Corporation Hauppauge New York", "December 28 1995", "Noncompliance With Performance Standards For Laser Products In That The Operator's Manuals Lacked Adequate Calibration Procedures And The Devices Had Several Labeling Noncompliances", "213 Units Were Distributed", "Nationwide Germany Korea"),
    @($null, $null, "II", $null, "Tempo Tingle Timers Installed With Diagnostic X Ray Systems Used To Regulate The Duration And To Some Extent The Amount Of Exposure To X Radiation", "Manufacturer Fda Approved The Firm's Corrective Action Plan January 3", "Tingle X Ray Products Inc", "January 3 1996", "Devices Were Not Tested For Accuracy In Accordance With The Specifications And Did Not Bear Proper Certification And Identification Labels As Required By 21 Cfr 1020", "17 Timers", $null),
    @($null, $null, "II", $null, "Immunocard Helicobacter Pylori Test Kit Catalog 710030 A Rapid Enzyme Immunoassay For The Detection Of Igg Antibodies To Helicobacter Pylori In Human Serum And Plasma", $null, "Meridian Diagnostics Inc", "November 14 1995", "The Enzyme Conjugate In The Kit Was Prepared At The Wrong Dilution Resulting In An Improper Sensitivity Of The Tests", "175 Test Kits Of Lot 710030 007 And 57 Test Kits Of Lot 710030 008 Were Distributed", "Nationwide Italy Russia Puerto Rico"),
    @($null, $null, "II", $null, "Captia Rubella M Eia Test Kit Product 801 165", $null, "Centocor Uk Ltd", "June 23 1995", "There Is A Higher Incidence Of Equivocal Positive Results Than Anticipated", "200 Kits Were Distributed", $null),
    @($null, $null, "II", $null, "Critikon Brand Protectiv And Protectiv Plus I", $null, "Johnson And Johnson Company Formerly Known As Critikon Inc", "April 10 1995", "An Opening May Be Present In The Clear Plastic Portion Of The Blister Package Therefore Compromising The Sterile Barrier Of The Package", "778 191 Catheters Were Distributed", $null),
    @($null, $null, "III", $null, "Cellfree Interleukin 2 Receptor Il 2 R Bead Assay Kit Catalog Ak 3120 For The Quantitative Measurement Of Interleukin 2 Receptor Level In Human Serum", $null, "T Cell Diagnostics Inc", "January 28 1994", "Controls Provided With The Product Were Not Meeting The Assigned Ranges And Fell Below Specified Limits", "20 Kits Were Distributed", "Massachusetts Texas Connecticut Louisiana Greece"),
    @($null, $null, "III", $null, "Glucometer Encore Test Strips Used To Measure Glucose In Whole Blood", $null, "Bayer Corporation Elkhart Indiana", "July 7 1995", "Subject Lots May Fail To Give Blood Glucose Reading And The Glucometer Encore Meter Will Display An E 3 Error Code", "45 451 Cartons", $null)
)

# A plain while loop kept losing track of the row/column counters in this
# runtime (and straight-line recursion over all 90 edits blows the host's
# expression-nesting budget), so the grid is walked as row-recursion
# nested around column-recursion -- lodash-style, each chain only as deep
# as the grid itself (<=13 / <=11).
#
# NOTE: recursive call arguments must be plain variables -- this host
# mis-parses `Func $var (expr)` as an indexing/invocation on $var, so every
# computed next-index is materialized into a temp variable first.
function Write-ReportCell($rowIndex, $colIndex) {
    $row = $reportGrid[$rowIndex]
    if ($colIndex -ge $row.Count) {
        return
    }
    $val = $row[$colIndex]
    if ($val -ne $null) {
        $ws.Cells.Item($rowIndex + 1, $colIndex + 1).Value = $val
    }
    $nextCol = $colIndex + 1
    Write-ReportCell $rowIndex $nextCol
}

function Write-ReportRow($rowIndex) {
    if ($rowIndex -ge $reportGrid.Count) {
        return
    }
    Write-ReportCell $rowIndex 0
    $nextRow = $rowIndex + 1
    Write-ReportRow $nextRow
}

Write-ReportRow 0

